# "clean up and add src" - record the missing "AL" (Urlaub) hours for
# Dezember and refresh the selection to point at the updated total cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the missing Dezember "AL" hours in B14 -----------------------
# B8 already carries the highlighted "value entered" look used for every
# other filled-in hours cell in this column; copy that formatting onto
# B14 before writing the new figure so it matches the rest of the sheet.
$ws.Range("B8").Copy()
$ws.Range("B14").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B14").Value = 27.67

# B19 (= B4-SUM(B5:B14)) recalculates automatically to 85.18 now that
# B14 is populated.

# --- Update the saved selection ----------------------------------------
[void]$ws.Range("B19").Select()
